$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2403945.8
$ws.Range("I9").Value = 749
$ws.Range("J9").Value = 4006077
$ws.Range("K9").Value = 749
$ws.Range("L9").Value = 4006077
$ws.Range("M9").Value = -580
$ws.Range("N9").Value = -4006415
$ws.Range("H17").Value = 5579.5835
$ws.Range("J17").Value = 5579.5835
$ws.Range("L17").Value = 16738.7505
$ws.Range("N17").Value = -17074.7505
$ws.Range("H111").Value = 29089.572
$ws.Range("I111").Value = 22725.4
$ws.Range("K111").Value = 68176.20000000001
$ws.Range("M111").Value = -65109.20000000001
$ws.Range("H116").Value = 4349.6924
$ws.Range("I116").Value = 3843.375
$ws.Range("J116").Value = 5159.8
$ws.Range("K116").Value = 3843.375
$ws.Range("L116").Value = 5159.8
$ws.Range("M116").Value = -401.375
$ws.Range("N116").Value = -12043.8
$ws.Range("H125").Value = 8359.4
$ws.Range("I125").Value = 7999.5
$ws.Range("J125").Value = 8599.333000000001
$ws.Range("K125").Value = 71995.5
$ws.Range("L125").Value = 77393.997
$ws.Range("M125").Value = -69535.5
$ws.Range("N125").Value = -82313.997
$ws.Range("H137").Value = 5964.8066
$ws.Range("I137").Value = 2886.7778
$ws.Range("J137").Value = 7224
$ws.Range("K137").Value = 8660.3334
$ws.Range("L137").Value = 21672
$ws.Range("M137").Value = -6110.3334
$ws.Range("N137").Value = -26772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 22000
$ws.Range("J23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("N23").Value = -22518
$ws.Range("H32").Value = 22489.607
$ws.Range("I32").Value = 22491.375
$ws.Range("K32").Value = 22491.375
$ws.Range("M32").Value = -22204.375
$ws.Range("H45").Value = 3333
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623
$ws.Range("H74").Value = 1001527.4
$ws.Range("I74").Value = 1001527.4
$ws.Range("K74").Value = 1001527.4
$ws.Range("M74").Value = -1000653.4
$ws.Range("H77").Value = 1001527.4
$ws.Range("I77").Value = 1001527.4
$ws.Range("K77").Value = 5007637
$ws.Range("M77").Value = -5003269
$ws.Range("H88").Value = 1883.8572
$ws.Range("J88").Value = 1696.75
$ws.Range("L88").Value = 1696.75
$ws.Range("N88").Value = -2508.75
$ws.Range("H91").Value = 1883.8572
$ws.Range("J91").Value = 1696.75
$ws.Range("L91").Value = 1696.75
$ws.Range("N91").Value = -4504.75
$ws.Range("H132").Value = 4702.1567
$ws.Range("I132").Value = 2869.125
$ws.Range("J132").Value = 11367.728
$ws.Range("K132").Value = 8607.375
$ws.Range("L132").Value = 34103.18399999999
$ws.Range("M132").Value = -6077.375
$ws.Range("N132").Value = -39163.18399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""
$ws.Range("H20").Value = 31246.457
$ws.Range("I20").Value = 45914.176
$ws.Range("J20").Value = 3133.3333
$ws.Range("K20").Value = 45914.176
$ws.Range("L20").Value = 3133.3333
$ws.Range("M20").Value = -45667.176
$ws.Range("N20").Value = -3627.3333
$ws.Range("H82").Value = 19049.053
$ws.Range("I82").Value = 19049.053
$ws.Range("K82").Value = 19049.053
$ws.Range("M82").Value = -18666.053
$ws.Range("H85").Value = 19049.053
$ws.Range("I85").Value = 19049.053
$ws.Range("K85").Value = 19049.053
$ws.Range("M85").Value = -17723.053
$ws.Range("H94").Value = 715.4545000000001
$ws.Range("I94").Value = 725.2381
$ws.Range("K94").Value = 725.2381
$ws.Range("M94").Value = -274.2381
$ws.Range("H97").Value = 6028.8335
$ws.Range("I97").Value = 6148.6
$ws.Range("J97").Value = 5430
$ws.Range("K97").Value = 6148.6
$ws.Range("L97").Value = 5430
$ws.Range("M97").Value = -5157.6
$ws.Range("N97").Value = -7412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16671914
$ws.Range("I31").Value = 43480760
$ws.Range("J31").Value = 6956.5674
$ws.Range("K31").Value = 43480760
$ws.Range("L31").Value = 6956.5674
$ws.Range("M31").Value = -43480465
$ws.Range("N31").Value = -7546.5674
$ws.Range("H34").Value = 16671914
$ws.Range("I34").Value = 43480760
$ws.Range("J34").Value = 6956.5674
$ws.Range("K34").Value = 43480760
$ws.Range("L34").Value = 6956.5674
$ws.Range("M34").Value = -43480558
$ws.Range("N34").Value = -7360.5674
$ws.Range("H59").Value = 22814.555
$ws.Range("I59").Value = 19871.572
$ws.Range("J59").Value = 33115
$ws.Range("K59").Value = 19871.572
$ws.Range("L59").Value = 33115
$ws.Range("M59").Value = -18726.572
$ws.Range("N59").Value = -35405
$ws.Range("H99").Value = 4833
$ws.Range("I99").Value = 4500
$ws.Range("J99").Value = 4999.5
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 4999.5
$ws.Range("M99").Value = -3002
$ws.Range("N99").Value = -7995.5
$ws.Range("H117").Value = 91500
$ws.Range("J117").Value = 91500
$ws.Range("L117").Value = 91500
$ws.Range("N117").Value = -100678
$ws.Range("H126").Value = 4833
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -19938.5
$ws.Range("H134").Value = 5889.793
$ws.Range("I134").Value = 5294
$ws.Range("K134").Value = 15882
$ws.Range("M134").Value = -13347

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1200
$ws.Range("I48").Value = 1200
$ws.Range("K48").Value = 3600
$ws.Range("M48").Value = -3350
$ws.Range("H61").Value = 86.28570999999999
$ws.Range("I61").Value = 22.555555
$ws.Range("K61").Value = 67.66666499999999
$ws.Range("M61").Value = 147.333335
$ws.Range("H68").Value = 2856.25
$ws.Range("J68").Value = 3483.3333
$ws.Range("L68").Value = 10449.9999
$ws.Range("N68").Value = -12071.9999
$ws.Range("H71").Value = 2856.25
$ws.Range("J71").Value = 3483.3333
$ws.Range("L71").Value = 31349.9997
$ws.Range("N71").Value = -39461.9997
$ws.Range("H122").Value = 50004496
$ws.Range("J122").Value = 1028.4117
$ws.Range("L122").Value = 9255.705300000001
$ws.Range("N122").Value = -14155.7053
$ws.Range("H129").Value = 29413682
$ws.Range("I129").Value = 600.2727
$ws.Range("K129").Value = 1800.8181
$ws.Range("M129").Value = 3199.1819
$ws.Range("H131").Value = 16678404
$ws.Range("J131").Value = 16231.429
$ws.Range("L131").Value = 48694.287
$ws.Range("N131").Value = -58774.287
$ws.Range("H139").Value = 2797.077
$ws.Range("J139").Value = 4006.6
$ws.Range("L139").Value = 12019.8
$ws.Range("N139").Value = -22299.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9999
$ws.Range("H73").Value = 9999
$ws.Range("H80").Value = 6512.75
$ws.Range("I80").Value = 2499.5
$ws.Range("J80").Value = 10526
$ws.Range("K80").Value = 2499.5
$ws.Range("L80").Value = 10526
$ws.Range("M80").Value = -1501.5
$ws.Range("N80").Value = -12522
$ws.Range("H83").Value = 6512.75
$ws.Range("I83").Value = 2499.5
$ws.Range("J83").Value = 10526
$ws.Range("K83").Value = 12497.5
$ws.Range("L83").Value = 52630
$ws.Range("M83").Value = -7505.5
$ws.Range("N83").Value = -62614
$ws.Range("H97").Value = 1148.5333
$ws.Range("I97").Value = 1095.5
$ws.Range("K97").Value = 1095.5
$ws.Range("M97").Value = -599.5
$ws.Range("H99").Value = 12610.182
$ws.Range("I99").Value = 13190.223
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 13190.223
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -10944.223
$ws.Range("N99").Value = -14492
$ws.Range("H113").Value = 6406.091
$ws.Range("I113").Value = 4850.5713
$ws.Range("J113").Value = 9128.25
$ws.Range("K113").Value = 4850.5713
$ws.Range("L113").Value = 9128.25
$ws.Range("M113").Value = -2680.5713
$ws.Range("N113").Value = -13468.25
$ws.Range("H126").Value = 4450.923
$ws.Range("I126").Value = 3559
$ws.Range("J126").Value = 4847.3335
$ws.Range("K126").Value = 10677
$ws.Range("L126").Value = 14542.0005
$ws.Range("M126").Value = -8207
$ws.Range("N126").Value = -19482.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1712.591
$ws.Range("I16").Value = 1149.2222
$ws.Range("J16").Value = 4247.75
$ws.Range("K16").Value = 1149.2222
$ws.Range("L16").Value = 4247.75
$ws.Range("M16").Value = -979.2221999999999
$ws.Range("N16").Value = -4587.75
$ws.Range("H61").Value = 7934.4
$ws.Range("I61").Value = 7934.4
$ws.Range("K61").Value = 7934.4
$ws.Range("M61").Value = -7732.4
$ws.Range("H93").Value = 2434.7827
$ws.Range("I93").Value = 2428.111
$ws.Range("J93").Value = 2458.8
$ws.Range("K93").Value = 2428.111
$ws.Range("L93").Value = 2458.8
$ws.Range("M93").Value = -1180.111
$ws.Range("N93").Value = -4954.8
$ws.Range("H113").Value = 7934.4
$ws.Range("I113").Value = 7934.4
$ws.Range("K113").Value = 7934.4
$ws.Range("M113").Value = -5764.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 5297.2
$ws.Range("J26").Value = 5297.2
$ws.Range("L26").Value = 5297.2
$ws.Range("N26").Value = -5883.2
$ws.Range("H96").Value = 4157.769
$ws.Range("I96").Value = 3519
$ws.Range("K96").Value = 3519
$ws.Range("M96").Value = -2146
$ws.Range("H122").Value = 4469.36
$ws.Range("I122").Value = 4538.9165
$ws.Range("K122").Value = 13616.7495
$ws.Range("M122").Value = -11166.7495
